$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A287").Value2 = 286
$ws.Range("B287").Value2 = 'Pós OS'
$ws.Range("C287").Value2 = 8005274938
$ws.Range("D287").Value2 = 46052.493726851862
$ws.Range("E287").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F287").Value2 = 'Detratores'
$ws.Range("G287").Value2 = 'a visita foi boa mas tenho tentado contato com voces porque a tecnica recomendou a substituicao do filtro e nao consigo retorno'
$ws.Range("H287").Value2 = 'Qualidade do Produto'
$ws.Range("I287").Value2 = 'Funcionamento geral'

$ws.Range("A288").Value2 = 287
$ws.Range("B288").Value2 = 'Pós OS'
$ws.Range("C288").Value2 = 8005266618
$ws.Range("D288").Value2 = 46052.501435185193
$ws.Range("E288").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F288").Value2 = 'Detratores'
$ws.Range("G288").Value2 = 'Marcam visita e não aparecem. Vc pede para trocar o filtro e precisa pagar! Estou com defeito no meu aparelho e até agora não veio ninguém! Agendei por duas vezes'
$ws.Range("H288").Value2 = 'Campo'
$ws.Range("I288").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A289").Value2 = 288
$ws.Range("B289").Value2 = 'Pós OS'
$ws.Range("C289").Value2 = 8005265640
$ws.Range("D289").Value2 = 46052.510208333333
$ws.Range("E289").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F289").Value2 = 'Neutros'
$ws.Range("G289").Value2 = 'Costuma dar algumas falhas'
$ws.Range("H289").Value2 = 'Qualidade do Produto'
$ws.Range("I289").Value2 = 'Funcionamento geral'

$ws.Range("A290").Value2 = 289
$ws.Range("B290").Value2 = 'Pós OS'
$ws.Range("C290").Value2 = 8005260783
$ws.Range("D290").Value2 = 46052.518680555557
$ws.Range("E290").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F290").Value2 = 'Detratores'
$ws.Range("G290").Value2 = 'Estamos aguardando a troca do aparelho pois todo o tempo que estamos com ele ele não funciona, já pedimos a troca do aparelho em maio de 2025 e não foi trocado e agora pedimos a subustituicap do aparelho para o convencional e estamos aguardando essa troca'
$ws.Range("H290").Value2 = 'Qualidade do Produto'
$ws.Range("I290").Value2 = 'Funcionamento geral'

$ws.Range("A291").Value2 = 290
$ws.Range("B291").Value2 = 'Pós OS'
$ws.Range("C291").Value2 = 8005259067
$ws.Range("D291").Value2 = 46052.547754629632
$ws.Range("E291").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F291").Value2 = 'Neutros'
$ws.Range("G291").Value2 = 'É satisfatório a qualidade do purificador e o serviço de assistência técnica. O valor da mensalidade é meio desestimulante.'
$ws.Range("H291").Value2 = 'Outros'
$ws.Range("I291").Value2 = 'Preço elevado'

$ws.Range("A292").Value2 = 291
$ws.Range("B292").Value2 = 'Pós OS'
$ws.Range("C292").Value2 = 8005265612
$ws.Range("D292").Value2 = 46052.554444444453
$ws.Range("E292").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F292").Value2 = 'Detratores'
$ws.Range("G292").Value2 = 'Excelente, tenho há muitos anos. Porém agora o agendamento está demorando muito'
$ws.Range("H292").Value2 = 'Capacidade'
$ws.Range("I292").Value2 = 'Data Distante'

$ws.Range("A293").Value2 = 292
$ws.Range("B293").Value2 = 'Pós OS'
$ws.Range("C293").Value2 = 8005283885
$ws.Range("D293").Value2 = 46052.573472222219
$ws.Range("E293").Value2 = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Range("F293").Value2 = 'Detratores'
$ws.Range("G293").Value2 = 'Depois da visita no dia 05.01.2026 o filtro apresentou um problema no dia 26.' + [char]10 + 'O tecnico que nos visitou informou que nenhuma limpeza havia sido feita ate entao.' + [char]10 + 'Apesar de eu ja ter reportado redução do fluxo de agua ninguem fez nada.' + [char]10 + 'O tecnico trocou o filtro e mais alguns sensores eletronicos, e o filtro voltou a funcionar da mesma maneira, ou seja ainda esta com o fluxo baixo.' + [char]10 + 'Infelizmente deixou muito a desejar.'
$ws.Range("H293").Value2 = 'Qualidade do Produto'
$ws.Range("I293").Value2 = 'Vazão/pressão da água'

$ws.Range("A294").Value2 = 293
$ws.Range("B294").Value2 = 'Pós OS'
$ws.Range("C294").Value2 = 8005274068
$ws.Range("D294").Value2 = 46052.574293981481
$ws.Range("E294").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F294").Value2 = 'Detratores'
$ws.Range("G294").Value2 = 'O tecnico nao apareceu'
$ws.Range("H294").Value2 = 'Campo'
$ws.Range("I294").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A295").Value2 = 294
$ws.Range("B295").Value2 = 'Pós OS'
$ws.Range("C295").Value2 = 8005283397
$ws.Range("D295").Value2 = 46052.64603009259
$ws.Range("E295").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F295").Value2 = 'Neutros'
$ws.Range("G295").Value2 = 'O único problema que tenho para relatar é o horário de chegada do técnico. Conforme agendado, estávamos esperando o técnico a partir das 13 h. Ele veio antes, por volta da 12 h que é o pico de atendimento da cozinha de nossa escola. Gostaríamos que fosse no horário agendado, digo, das 13 h às 18 h. Do resto tudo perfeito !'
$ws.Range("H295").Value2 = 'Campo'
$ws.Range("I295").Value2 = 'Fora do período agendado'

$ws.Range("A296").Value2 = 295
$ws.Range("B296").Value2 = 'Pós OS'
$ws.Range("C296").Value2 = 8005278837
$ws.Range("D296").Value2 = 46052.827280092592
$ws.Range("E296").Value2 = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Range("F296").Value2 = 'Neutros'
$ws.Range("G296").Value2 = 'Experiência boa e satisfatória, porém o preço está muito alto.'
$ws.Range("H296").Value2 = 'Outros'
$ws.Range("I296").Value2 = 'Preço elevado'

$ws.Range("A297").Value2 = 296
$ws.Range("B297").Value2 = 'Pós OS'
$ws.Range("C297").Value2 = 8005279298
$ws.Range("D297").Value2 = 46053.400891203702
$ws.Range("E297").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F297").Value2 = 'Detratores'
$ws.Range("G297").Value2 = 'A pessoa não veio.'
$ws.Range("H297").Value2 = 'Campo'
$ws.Range("I297").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A298").Value2 = 297
$ws.Range("B298").Value2 = 'Pós OS'
$ws.Range("C298").Value2 = 8005265467
$ws.Range("D298").Value2 = 46053.435358796298
$ws.Range("E298").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F298").Value2 = 'Neutros'
$ws.Range("G298").Value2 = 'boa'
$ws.Range("H298").Value2 = 'Outros'
$ws.Range("I298").Value2 = 'Satisfação geral'

$ws.Range("A299").Value2 = 298
$ws.Range("B299").Value2 = 'Pós OS'
$ws.Range("C299").Value2 = 8005245207
$ws.Range("D299").Value2 = 46053.436574074083
$ws.Range("E299").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F299").Value2 = 'Neutros'
$ws.Range("G299").Value2 = 'Da outra vez marcaram e ninguem apareceu e nem deu satisfaçao' + [char]10 + 'Desta vez foi tudo bem'
$ws.Range("H299").Value2 = 'Campo'
$ws.Range("I299").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A300").Value2 = 299
$ws.Range("B300").Value2 = 'Pós OS'
$ws.Range("C300").Value2 = 8005253746
$ws.Range("D300").Value2 = 46053.854456018518
$ws.Range("E300").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F300").Value2 = 'Detratores'
$ws.Range("G300").Value2 = 'Depois de um mês de espera, foi agendada a troca do purificador que não ocorreu. Agendada a troca pela segunda vez, foi instalado o purificador, porém a água está com gosto estranho e o controle de volume de água não está funcionando. Quero cancelar minha assinatura do purificador.'
$ws.Range("H300").Value2 = 'Supply'
$ws.Range("I300").Value2 = 'Falta de estoque/peças'

$ws.Range("A301").Value2 = 300
$ws.Range("B301").Value2 = 'Instalação'
$ws.Range("C301").Value2 = 8005256304
$ws.Range("D301").Value2 = 46054.640370370369
$ws.Range("E301").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F301").Value2 = 'Neutros'
$ws.Range("G301").Value2 = 'Ultimamente a Brastemp não comparece na data agendada para a manutenção da limpeza do filtro.' + [char]10 + 'Fico aguardando e não recebo nenhuma satisfação.' + [char]10 + 'É preciso reagendar e espero mais um mês.' + [char]10 + 'Sou profissional autônoma, o que me leva a cancelar meus compromissos no período que fico aguardando a Brastemp.' + [char]10 + 'Tenho prejuízo financeiro.'
$ws.Range("H301").Value2 = 'Campo'
$ws.Range("I301").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A302").Value2 = 301
$ws.Range("B302").Value2 = 'Pós OS'
$ws.Range("C302").Value2 = 8005280499
$ws.Range("D302").Value2 = 46054.72960648148
$ws.Range("E302").Value2 = 'AT_ECO_CE_FORTALEZA'
$ws.Range("F302").Value2 = 'Detratores'
$ws.Range("G302").Value2 = 'Demorou 4 dias para consertar e eu pagando água mineral'
$ws.Range("H302").Value2 = 'Campo'
$ws.Range("I302").Value2 = 'Reincidência'

$ws.Range("A303").Value2 = 302
$ws.Range("B303").Value2 = 'Pós OS'
$ws.Range("C303").Value2 = 8005282675
$ws.Range("D303").Value2 = 46054.959293981483
$ws.Range("E303").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F303").Value2 = 'Detratores'
$ws.Range("G303").Value2 = 'Mais uma vez o agendamento não foi cumprido.'
$ws.Range("H303").Value2 = 'Campo'
$ws.Range("I303").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A304").Value2 = 303
$ws.Range("B304").Value2 = 'Pós OS'
$ws.Range("C304").Value2 = 8005277668
$ws.Range("D304").Value2 = 46055.378750000003
$ws.Range("E304").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F304").Value2 = 'Detratores'
$ws.Range("G304").Value2 = 'Tecnico beinem casa e não resolveu o problema.'
$ws.Range("H304").Value2 = 'Campo'
$ws.Range("I304").Value2 = 'Reincidência'

$ws.Range("A305").Value2 = 304
$ws.Range("B305").Value2 = 'Pós OS'
$ws.Range("C305").Value2 = 8005258927
$ws.Range("D305").Value2 = 46055.390648148154
$ws.Range("E305").Value2 = 'FRQ_ECO_RJ_ATLANTICA'
$ws.Range("F305").Value2 = 'Detratores'
$ws.Range("G305").Value2 = 'Não tive visita nenhuma, estou a meses solicitando um técnico e nada!! Tive q fazer mudança de endereço e paguei particular a instalação, vcs não me deram nenhuma atenção. Várias atendentes ligam e não conseguem nem mudar meu endereço. Atendimento péssimo'
$ws.Range("H305").Value2 = 'Atendimento'
$ws.Range("I305").Value2 = 'Qualidade do atendimento'

$ws.Range("A306").Value2 = 305
$ws.Range("B306").Value2 = 'Pós OS'
$ws.Range("C306").Value2 = 8005211570
$ws.Range("D306").Value2 = 46055.400266203702
$ws.Range("E306").Value2 = 'FRQ_ECO_DF_BRASILIA_2'
$ws.Range("F306").Value2 = 'Neutros'
$ws.Range("G306").Value2 = 'Marcaram pro sábado de manhã e o técnico veio na sexta à tarde. Por sorte eu estava em casa.'
$ws.Range("H306").Value2 = 'Campo'
$ws.Range("I306").Value2 = 'Fora do período agendado'

$ws.Range("A307").Value2 = 306
$ws.Range("B307").Value2 = 'Pós OS'
$ws.Range("C307").Value2 = 8005285049
$ws.Range("D307").Value2 = 46055.502395833333
$ws.Range("E307").Value2 = 'FRQ_ECO_RJ_ATLANTICA'
$ws.Range("F307").Value2 = 'Detratores'
$ws.Range("G307").Value2 = 'O técnico nao apareceu bao deu satisfação e continuo sem água no filtro que está com defeito. Muita insatisfeito.'
$ws.Range("H307").Value2 = 'Campo'
$ws.Range("I307").Value2 = 'Técnico não cumpriu a agenda'

$ws.Range("A308").Value2 = 307
$ws.Range("B308").Value2 = 'Pós OS'
$ws.Range("C308").Value2 = 8005281526
$ws.Range("D308").Value2 = 46055.511180555557
$ws.Range("E308").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F308").Value2 = 'Detratores'
$ws.Range("G308").Value2 = 'O purificador não está funcionando e já faz um mês que não consigo resolver o problema.'
$ws.Range("H308").Value2 = 'Campo'
$ws.Range("I308").Value2 = 'Reincidência'

$ws.Range("A309").Value2 = 308
$ws.Range("B309").Value2 = 'Pós OS'
$ws.Range("C309").Value2 = 8005284507
$ws.Range("D309").Value2 = 46055.512384259258
$ws.Range("E309").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F309").Value2 = 'Detratores'
$ws.Range("G309").Value2 = 'Vcs estão fazendo a maior bagunça com essas marcações por secretária eletrônica.' + [char]10 + 'Chegou um rapaz aqui na semana passada na minha casa em Piratininga para fazer manutenção em meu filtro da casa.Nao deixei entrar porquê o que marquei foi para essa semana na próxima quinta feira no meu apartamento em Icaraí.' + [char]10 + 'Tenho 2 assinaturas e vcs fazem a maior confusão na hora de mandarem o técnico .' + [char]10 + 'Fico mal assistida nas duas'
$ws.Range("H309").Value2 = 'Campo'
$ws.Range("I309").Value2 = 'Fora do período agendado'

$ws.Range("A310").Value2 = 309
$ws.Range("B310").Value2 = 'Pós OS'
$ws.Range("C310").Value2 = 8005278395
$ws.Range("D310").Value2 = 46055.512511574067
$ws.Range("E310").Value2 = 'FRQ_ECO_RJ_ATLANTICA'
$ws.Range("F310").Value2 = 'Detratores'
$ws.Range("G310").Value2 = 'Infelizmente foi instalado com defeito sem funcionar o gelo liguei para consertar e mudar de lugar não pode fazer as duas solicitações então optei por mudar de lugar e depois vou pedir para consertar infelizmente muita burocracia'
$ws.Range("H310").Value2 = 'Qualidade do Produto'
$ws.Range("I310").Value2 = 'Funcionamento geral'

$ws.Range("A311").Value2 = 310
$ws.Range("B311").Value2 = 'Pós OS'
$ws.Range("C311").Value2 = 8005251148
$ws.Range("D311").Value2 = 46055.540162037039
$ws.Range("E311").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F311").Value2 = 'Detratores'
$ws.Range("G311").Value2 = 'Estão demorando muito para trocar o sensor que tem gerado o congelamento do filtro.' + [char]10 + 'Estou sem usá-lo a quase três semanas'
$ws.Range("H311").Value2 = 'Supply'
$ws.Range("I311").Value2 = 'Falta de estoque/peças'

$ws.Range("A312").Value2 = 311
$ws.Range("B312").Value2 = 'Pós OS'
$ws.Range("C312").Value2 = 8005278311
$ws.Range("D312").Value2 = 46055.592326388891
$ws.Range("E312").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F312").Value2 = 'Detratores'
$ws.Range("G312").Value2 = 'O Rapaz veio instalar o pressurizador e segundo ele veio a peça errada, ligou para uma pessoa para falar sobre, ficou de voltar no dia seguinte para instalar a peça correta e até hoje não retornou.'
$ws.Range("H312").Value2 = 'Campo'
$ws.Range("I312").Value2 = 'Reincidência'

$ws.Range("A313").Value2 = 312
$ws.Range("B313").Value2 = 'Pós OS'
$ws.Range("C313").Value2 = 8005275635
$ws.Range("D313").Value2 = 46055.703217592592
$ws.Range("E313").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F313").Value2 = 'Detratores'
$ws.Range("G313").Value2 = 'O técnico solicitou a troca do aparelho. Pois o mesmo ainda não se encontra funcionando.'
$ws.Range("H313").Value2 = 'Supply'
$ws.Range("I313").Value2 = 'Falta de estoque/peças'

$ws.Range("A314").Value2 = 313
$ws.Range("B314").Value2 = 'Instalação'
$ws.Range("C314").Value2 = 8005280530
$ws.Range("D314").Value2 = 46055.800706018519
$ws.Range("E314").Value2 = 'FRQ_ECO_SP_ZONA_SUL_03'
$ws.Range("F314").Value2 = 'Detratores'
$ws.Range("G314").Value2 = 'Não posso utilizar. A água tem um gosto de óleo químico'
$ws.Range("H314").Value2 = 'Qualidade do Produto'
$ws.Range("I314").Value2 = 'Desempenho da filtragem'

$ws.Range("A315").Value2 = 314
$ws.Range("B315").Value2 = 'Pós OS'
$ws.Range("C315").Value2 = 8005272119
$ws.Range("D315").Value2 = 46055.870787037027
$ws.Range("E315").Value2 = 'FRQ_ECO_BA_SALVADOR'
$ws.Range("F315").Value2 = 'Neutros'
$ws.Range("G315").Value2 = 'Tenho o purificador desde 2014 ,houve um problema no filtro , pois este perdeu a validade fazendo com que a água ficasse com mal cheiro e gosto desagradável. Diante do exposto solicitei um técnico com urgência pois tenho criança mas infelizmente levei mais de 15 dias para um técnico trocasse o mesmo.'
$ws.Range("H315").Value2 = 'Capacidade'
$ws.Range("I315").Value2 = 'Data Distante'

$ws.Range("A316").Value2 = 315
$ws.Range("B316").Value2 = 'Pós OS'
$ws.Range("C316").Value2 = 8005282966
$ws.Range("D316").Value2 = 46056.384085648147
$ws.Range("E316").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F316").Value2 = 'Detratores'
$ws.Range("G316").Value2 = 'Péssimo atendimento, estou extremamente decepcionado, sou cliente a mais de 15 anos e estou tendo uma dificuldade enorme para ter o meu filtro operacional novamente. Primeiramente a demora em conseguir uma visita técnica, mesmo assim aguardamos, porém o técnico veio, não deu o laudo condenando o filtro (não está gelando a água a mais de duas semanas), o técnico efetuou alguma manutenção mal feita e foi embora do local dizendo que em 20 minutos estaria gelando, porém isso não aconteceu, o filtro continua sem gelar. Pagamos por um serviço que não está sendo entregue. Já ligamos novamente e só tem data para uma nova visita para daqui a uma semana!!!'
$ws.Range("H316").Value2 = 'Qualidade do Produto'
$ws.Range("I316").Value2 = 'Não gela'

$ws.Range("A317").Value2 = 316
$ws.Range("B317").Value2 = 'Pós OS'
$ws.Range("C317").Value2 = 8005254036
$ws.Range("D317").Value2 = 46056.384328703702
$ws.Range("E317").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F317").Value2 = 'Detratores'
$ws.Range("G317").Value2 = 'Não respondem e-mails, agendamentos errados.' + [char]10 + 'Estamos com equipamento parado a meses.  Solicitamos o contato do nosso gerente de conta e não tivemos retorno,  gostaríamos de mais informações sobre o equipamento de água gaseificada.  Mas a Brastemp parece que não quer vender serviços'
$ws.Range("H317").Value2 = 'Atendimento'
$ws.Range("I317").Value2 = 'Erro de comunicação'

$ws.Range("A318").Value2 = 317
$ws.Range("B318").Value2 = 'Pós OS'
$ws.Range("C318").Value2 = 8005253400
$ws.Range("D318").Value2 = 46056.391793981478
$ws.Range("E318").Value2 = 'FRQ_ECO_PE_RECIFE'
$ws.Range("F318").Value2 = 'Detratores'
$ws.Range("G318").Value2 = 'Precisei ligar para fazer agendamento sendo que cada semestre a própria Brastemp marcava o agendamento'
$ws.Range("H318").Value2 = 'Atendimento'
$ws.Range("I318").Value2 = 'Qualidade do atendimento'

$ws.Range("A319").Value2 = 318
$ws.Range("B319").Value2 = 'Pós OS'
$ws.Range("C319").Value2 = 8005277186
$ws.Range("D319").Value2 = 46056.411921296298
$ws.Range("E319").Value2 = 'FRQ_ECO_RJ_OESTE'
$ws.Range("F319").Value2 = 'Detratores'
$ws.Range("G319").Value2 = 'a demora mais de um mês para o técnico trocar um filtro.' + [char]10 + 'sem cabimento não demorou 10 minutos.' + [char]10 + 'péssimo atendimento.'
$ws.Range("H319").Value2 = 'Campo'
$ws.Range("I319").Value2 = 'Qualidade da manutenção'

$ws.Range("A320").Value2 = 319
$ws.Range("B320").Value2 = 'Pós OS'
$ws.Range("C320").Value2 = 8005262693
$ws.Range("D320").Value2 = 46056.420810185176
$ws.Range("E320").Value2 = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Range("F320").Value2 = 'Detratores'
$ws.Range("G320").Value2 = 'Gosto muito do filtro da Brastemp, mas ultimamente as visitas emergenciais por algum problema  no aparelho demoram muito. Cheguei a ficar 1 mês sem filtro, o mesmo estava com a água quente, ontem 02/02 o técnico foi muito atencioso e trocou o sensor, gostaria de uma atendimento mais eficaz porque pago mensalmente sempre sem dever nada a empresa, gostaria de ser tratada da mesma forma.'
$ws.Range("H320").Value2 = 'Supply'
$ws.Range("I320").Value2 = 'Falta de estoque/peças'

$ws.Range("A321").Value2 = 320
$ws.Range("B321").Value2 = 'Pós OS'
$ws.Range("C321").Value2 = 8005247370
$ws.Range("D321").Value2 = 46056.456076388888
$ws.Range("E321").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F321").Value2 = 'Detratores'
$ws.Range("G321").Value2 = 'Sou cliente a mais de 10 anos e pasmem, até hj continuo c o mesmo aparelho, que já apresentou problemas de aquecimento varias vezes e a desculpa a cd manutenção é a mesma, troca-se uma peça e aguarde pois estaremos solicitando a troca do aparelho, até pq em 10 anos já sairam inúmeras versões certamente + modernas, eficazes ... já avisei que irei fz o cancelamento, o antepenúltimo técnico que esteve fz a manutenção ficou p um bom tempo em ligação confirmando que seria trocado o aparelho e passado 1 ano, continuo na mesma ...o filtro foi trocado pouquissimas vezes e qual garantia temos que está sendo benéfico? o teste do copo???? caso n seja trocado o aparelho, estarei procedendo ao cancelamento.'
$ws.Range("H321").Value2 = 'Qualidade do Produto'
$ws.Range("I321").Value2 = 'Não gela'

$ws.Range("A322").Value2 = 321
$ws.Range("B322").Value2 = 'Pós OS'
$ws.Range("C322").Value2 = 8005282361
$ws.Range("D322").Value2 = 46056.502210648148
$ws.Range("E322").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F322").Value2 = 'Detratores'
$ws.Range("G322").Value2 = 'O técnico fez duas visitas e não resolveu o problema'
$ws.Range("H322").Value2 = 'Campo'
$ws.Range("I322").Value2 = 'Reincidência'

$ws.Range("A323").Value2 = 322
$ws.Range("B323").Value2 = 'Pós OS'
$ws.Range("C323").Value2 = 8005287614
$ws.Range("D323").Value2 = 46056.504351851851
$ws.Range("E323").Value2 = 'FRQ_ECO_SP_OSASCO'
$ws.Range("F323").Value2 = 'Neutros'
$ws.Range("G323").Value2 = 'O único senão é  que quando falta luz ele não funciona. A manutenção dele foi negligenciada pois desde 2023 não  recebi aviso que estava na hora da manutenção.  Vieram agora para cessar o vazamento mas a vistoria,  manutenção  nada foi feito.'
$ws.Range("H323").Value2 = 'Qualidade do Produto'
$ws.Range("I323").Value2 = 'Vazamento'

$ws.Range("A324").Value2 = 323
$ws.Range("B324").Value2 = 'Pós OS'
$ws.Range("C324").Value2 = 8005285952
$ws.Range("D324").Value2 = 46056.510185185187
$ws.Range("E324").Value2 = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Range("F324").Value2 = 'Neutros'
$ws.Range("G324").Value2 = 'Foi muito bom ,só que marquei a tarde e foi cedo'
$ws.Range("H324").Value2 = 'Campo'
$ws.Range("I324").Value2 = 'Fora do período agendado'

$ws.Range("A325").Value2 = 324
$ws.Range("B325").Value2 = 'Pós OS'
$ws.Range("C325").Value2 = 8005275703
$ws.Range("D325").Value2 = 46056.579317129632
$ws.Range("E325").Value2 = 'FRQ_ECO_SP_S B CAMPO'
$ws.Range("F325").Value2 = 'Detratores'
$ws.Range("G325").Value2 = 'Meu filtro continua sem funcionar desde o dia da instalação!!!! Solicito troca imediata do equipamento.'
$ws.Range("H325").Value2 = 'Campo'
$ws.Range("I325").Value2 = 'Reincidência'

$ws.Range("A326").Value2 = 325
$ws.Range("B326").Value2 = 'Pós OS'
$ws.Range("C326").Value2 = 8005285741
$ws.Range("D326").Value2 = 46056.581597222219
$ws.Range("E326").Value2 = 'FRQ_ECO_SP_CAMPINAS_2'
$ws.Range("F326").Value2 = 'Detratores'
$ws.Range("G326").Value2 = 'Tenho a assinatura ininterrupta há mais de 22 anos, ou seja, desde o início desse serviço aqui em Campinas.' + [char]10 + 'Esta última manutenção foi feita de maneira muito rápida, parecendo não se atentar para todos os itens protocolares.' + [char]10 + 'O preço da assinatura está muito alto, pelo que reinvidico desconto, tendo em vista a fidelização!'
$ws.Range("H326").Value2 = 'Campo'
$ws.Range("I326").Value2 = 'Qualidade da manutenção'

$ws.Range("A327").Value2 = 326
$ws.Range("B327").Value2 = 'Pós OS'
$ws.Range("C327").Value2 = 8005265966
$ws.Range("D327").Value2 = 46056.589722222219
$ws.Range("E327").Value2 = 'FRQ_ECO_SP_SANTOS'
$ws.Range("F327").Value2 = 'Neutros'
$ws.Range("G327").Value2 = 'Valor'
$ws.Range("H327").Value2 = 'Outros'
$ws.Range("I327").Value2 = 'Preço elevado'

$ws.Range("A328").Value2 = 327
$ws.Range("B328").Value2 = 'Pós OS'
$ws.Range("C328").Value2 = 8005255071
$ws.Range("D328").Value2 = 46056.61986111111
$ws.Range("E328").Value2 = 'FRQ_ECO_RJ_ATLANTICA'
$ws.Range("F328").Value2 = 'Neutros'
$ws.Range("G328").Value2 = 'A visita não resolveu meu problema, é preciso mudar o aparelho.'
$ws.Range("H328").Value2 = 'Campo'
$ws.Range("I328").Value2 = 'Reincidência'

# Update selection to match the final active range used during editing
$ws.Range("A287:I328").Select()
